# Rename the "Group" sheet to "Untagged" (part of adding a group-movement
# parser that now distinguishes tagged vs. untagged individuals), and make
# it the active/selected tab instead of "Individual".
$wb = $excel.ActiveWorkbook

$groupSheet = $wb.Worksheets.Item("Group")
$groupSheet.Name = "Untagged"

$groupSheet.Activate()
